$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("3. Apresenta preço do componente") is removed entirely; everything
# below it shifts up by one row, and the numbered step text is renumbered.
$ws.Rows.Item(9).Delete()

# Renumbered / reworded cells after the row shift (new row numbers).
$ws.Range("D9").Value  = "3. Verifica necessidade de componentes extras"
$ws.Range("D10").Value = "4. Confirma compatibilidade com componentes escolhidos"
$ws.Range("D11").Value = "5. Apresenta preço final"
$ws.Range("C12").Value = "6. Confirma componte"
$ws.Range("D13").Value = "7.  Adiciona componente"

$ws.Range("B15").Value = " Alternativa 1 [Componente incompatível com existente]       Passo 4"
$ws.Range("D15").Value = "4.1. Informa de incompatibilidade e apresenta os componentes incompatíveis"
$ws.Range("C16").Value = "4.2. Mantém componente atual"
$ws.Range("D17").Value = "4.3. Retira componentes incompatíveis"
$ws.Range("D18").Value = "4.4 Retorna ao passo 5"

$ws.Range("B19").Value = " Alternativa 2 [Necessita mais componentes] Passo 3"
$ws.Range("D19").Value = "3.1. Informa de necessidade de mais componentes e apresenta os mesmos"
$ws.Range("C20").Value = "3.2. Confirma componentes extras"
$ws.Range("D21").Value = "3.3. Adiciona componentes extra "

$ws.Range("B22").Value = " Excepção 3 [Cliente não aceita componente] Passos 3.2, 4.2 e 6"

# Restore the view state recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C28").Select()
